$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @("ECs", "Postn", "Ptk7", "ECs", 3, 1, 11.66506566666667, 34.995197, 0.002461364886242511, 0.002461364886242511, 3, 1, 1.488193, 4.464579, 0.09810006315031751, 0.09810006315031751, 17.35986906967367, 156.238821627063, 0.0002414600507763644, 0.0002414600507763644),
        @("ECs", "Postn", "Ptk7", "FAPs", 3, 1, 11.66506566666667, 34.995197, 0.002461364886242511, 0.002461364886242511, 3, 1, 9.650305666666668, 28.950917, 0.6361376483560043, 0.6361376483560043, 112.5714493050722, 1013.143043745649, 0.001565766870480355, 0.001565766870480355),
        @("ECs", "Postn", "Ptk7", "M2", 3, 1, 11.66506566666667, 34.995197, 0.002461364886242511, 0.002461364886242511, 1, 0.3333333333333333, 0.009611333333333333, 0.028834, 0.0006335686345512656, 0.0006335686345512656, 0.1121168344775556, 1.009051510298, [double]"1.559443590109099e-06", [double]"1.559443590109099e-06"),
        @("ECs", "Postn", "Ptk7", "sCs", 3, 1, 11.66506566666667, 34.995197, 0.002461364886242511, 0.002461364886242511, 3, 1, 4.022043333333333, 12.06613, 0.2651287198591268, 0.2651287198591268, 46.91739959751223, 422.25659637761, 0.0006525785213956823, 0.0006525785213956823),
        @("FAPs", "Postn", "Ptk7", "ECs", 3, 1, 4683.852620333333, 14051.557861, 0.9883073701819801, 0.9883073701819801, 3, 1, 1.488193, 4.464579, 0.09810006315031751, 0.09810006315031751, 6970.476682611724, 62734.29014350551, 0.09695301542677648, 0.09695301542677648),
        @("FAPs", "Postn", "Ptk7", "FAPs", 3, 1, 4683.852620333333, 14051.557861, 0.9883073701819801, 0.9883073701819801, 3, 1, 9.650305666666668, 28.950917, 0.6361376483560043, 0.6361376483560043, 45200.60948383428, 406805.4853545086, 0.6286995263204718, 0.6286995263204718),
        @("FAPs", "Postn", "Ptk7", "M2", 3, 1, 4683.852620333333, 14051.557861, 0.9883073701819801, 0.9883073701819801, 1, 0.3333333333333333, 0.009611333333333333, 0.028834, 0.0006335686345512656, 0.0006335686345512656, 45.01806881823044, 405.162619364074, 0.0006261605510431493, 0.0006261605510431493),
        @("FAPs", "Postn", "Ptk7", "sCs", 3, 1, 4683.852620333333, 14051.557861, 0.9883073701819801, 0.9883073701819801, 3, 1, 4.022043333333333, 12.06613, 0.2651287198591268, 0.2651287198591268, 18838.65820592755, 169547.9238533479, 0.2620286678836886, 0.2620286678836886),
        @("M2", "Postn", "Ptk7", "ECs", 3, 1, 3.125034333333333, 9.375103, 0.000659391896811063, 0.000659391896811063, 3, 1, 1.488193, 4.464579, 0.09810006315031751, 0.09810006315031751, 4.650654219626333, 41.85588797663699, [double]"6.468638671797292e-05", [double]"6.468638671797292e-05"),
        @("M2", "Postn", "Ptk7", "FAPs", 3, 1, 3.125034333333333, 9.375103, 0.000659391896811063, 0.000659391896811063, 3, 1, 9.650305666666668, 28.950917, 0.6361376483560043, 0.6361376483560043, 30.15753653549456, 271.417828819451, 0.0004194640105823946, 0.0004194640105823946),
        @("M2", "Postn", "Ptk7", "M2", 3, 1, 3.125034333333333, 9.375103, 0.000659391896811063, 0.000659391896811063, 1, 0.3333333333333333, 0.009611333333333333, 0.028834, 0.0006335686345512656, 0.0006335686345512656, 0.03003574665577777, 0.270321719902, [double]"4.177700236967542e-07", [double]"4.177700236967542e-07"),
        @("M2", "Postn", "Ptk7", "sCs", 3, 1, 3.125034333333333, 9.375103, 0.000659391896811063, 0.000659391896811063, 3, 1, 4.022043333333333, 12.06613, 0.2651287198591268, 0.2651287198591268, 12.56902350682111, 113.12121156139, 0.0001748237294869986, 0.0001748237294869986),
        @("sCs", "Postn", "Ptk7", "ECs", 3, 1, 40.624396, 121.873188, 0.008571873034966258, 0.00857187303496626, 3, 1, 1.488193, 4.464579, 0.09810006315031751, 0.09810006315031751, 60.45694175642799, 544.112475807852, 0.0008409012860466937, 0.0008409012860466939),
        @("sCs", "Postn", "Ptk7", "FAPs", 3, 1, 40.624396, 121.873188, 0.008571873034966258, 0.00857187303496626, 3, 1, 9.650305666666668, 28.950917, 0.6361376483560043, 0.6361376483560043, 392.0378389237107, 3528.340550313396, 0.005452891154469681, 0.005452891154469681),
        @("sCs", "Postn", "Ptk7", "M2", 3, 1, 40.624396, 121.873188, 0.008571873034966258, 0.00857187303496626, 1, 0.3333333333333333, 0.009611333333333333, 0.028834, 0.0006335686345512656, 0.0006335686345512656, 0.3904546114213333, 3.514091502792, [double]"5.430869894310385e-06", [double]"5.430869894310386e-06"),
        @("sCs", "Postn", "Ptk7", "sCs", 3, 1, 40.624396, 121.873188, 0.008571873034966258, 0.00857187303496626, 3, 1, 4.022043333333333, 12.06613, 0.2651287198591268, 0.2651287198591268, 163.3930811024933, 1470.53772992244, 0.002272649724555573, 0.002272649724555573)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}
